# Added 4wk low sales check: refresh forecast figures on "Forecast Comparison"
# and update the derived totals on "Summary".

$wb = $excel.ActiveWorkbook

# --- Forecast Comparison sheet -------------------------------------------
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Week W10 (row 2)
$ws.Range("D2").Value = 139
$ws.Range("H2").Value = 8.16
$ws.Range("L2").Value = 0.82

# Week W11 (row 3)
$ws.Range("D3").Value = 127
$ws.Range("H3").Value = 7.83
$ws.Range("L3").Value = 0.86

# Week W12 (row 4)
$ws.Range("D4").Value = 118
$ws.Range("H4").Value = 7.36
$ws.Range("L4").Value = 0.8

# Week W13 (row 5)
$ws.Range("D5").Value = 122
$ws.Range("H5").Value = 6.15
$ws.Range("L5").Value = 0.9

# Week W14 (row 6)
$ws.Range("D6").Value = 120
$ws.Range("H6").Value = 5.23
$ws.Range("L6").Value = 0.85

# Week W15 (row 7)
$ws.Range("D7").Value = 104
$ws.Range("H7").Value = 4.88
$ws.Range("L7").Value = 0.81

# Week W16 (row 8)
$ws.Range("D8").Value = 80
$ws.Range("H8").Value = 5.05
$ws.Range("L8").Value = 0.96

# Week W17 (row 9)
$ws.Range("D9").Value = 64
$ws.Range("H9").Value = 5.06
$ws.Range("L9").Value = 1.04

# Week W18 (row 10)
$ws.Range("H10").Value = 4.19
$ws.Range("L10").Value = 1.02

# Week W19 (row 11)
$ws.Range("D11").Value = 64
$ws.Range("H11").Value = 3.09
$ws.Range("L11").Value = 1.18

# Week W20 (row 12)
$ws.Range("D12").Value = 59
$ws.Range("H12").Value = 2.27
$ws.Range("L12").Value = 0.83

# Week W21 (row 13)
$ws.Range("D13").Value = 49
$ws.Range("H13").Value = 1.53
$ws.Range("L13").Value = 1.2

# Week W22 (row 14)
$ws.Range("D14").Value = 43
$ws.Range("H14").Value = 0.6
$ws.Range("J14").Value = "Urgent"

# Week W23 (row 15)
$ws.Range("D15").Value = 39
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = "High"
$ws.Range("J15").Value = "Urgent"
$ws.Range("L15").Value = 0.82

# Week W24 (row 16)
$ws.Range("D16").Value = 36
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = "High"
$ws.Range("J16").Value = "Urgent"
$ws.Range("L16").Value = 0.96

# Week W25 (row 17)
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = "High"
$ws.Range("J17").Value = "Urgent"
$ws.Range("L17").Value = 0.88

# --- Summary sheet ---------------------------------------------------------
# These cells hold numeric-looking text, so force text formatting while
# entering the value, then restore the default "Normal" style so no stray
# cell formatting is introduced.
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("B9").NumberFormat = "@"
$summary.Range("B9").Value = "1261"
$summary.Range("B9").Style = "Normal"

$summary.Range("B10").NumberFormat = "@"
$summary.Range("B10").Value = "874"
$summary.Range("B10").Style = "Normal"

$summary.Range("B11").NumberFormat = "@"
$summary.Range("B11").Value = "506"
$summary.Range("B11").Style = "Normal"

$summary.Range("B12").NumberFormat = "@"
$summary.Range("B12").Value = "139"
$summary.Range("B12").Style = "Normal"
